# "chnages for day 3"
# Adds a new "ValidLogin" worksheet (after the existing "tc1" sheet) that
# holds a second username/password pair, and makes that new sheet the
# active / selected tab.

$wb = $excel.ActiveWorkbook

# tc1 is the first (and, so far, only) sheet in the workbook.
$tc1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after tc1 so tab order stays [tc1, ValidLogin].
$validLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc1)
$validLogin.Name = "ValidLogin"

# Header row, matching tc1's layout.
$validLogin.Range("A1").Value = "username"
$validLogin.Range("B1").Value = "password"

# Credentials row. Write B2 before A2 so new shared-string entries land in
# the same order as the authored workbook ("pointofsale" then "ADMIN").
$validLogin.Range("B2").Value = "pointofsale"
$validLogin.Range("A2").Value = "ADMIN"

# Leave the cursor on A2 and zoom the new sheet in, then make it the
# active/visible tab.
$validLogin.Range("A2").Select()
$excel.ActiveWindow.Zoom = 205
$validLogin.Activate()
